$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column ("Date") currently holds the literal text "6-30-2007-08"
# for each data row (2 through 31). The stats were off by a day because of
# how the NBA stats date was originally captured, so normalize the values
# to the correct ISO-ish text "2008-06-30". The values must stay plain
# text (not get auto-converted to a date serial number), so the range is
# pre-formatted as Text before the write and the formatting is cleared
# again afterwards (restoring the default/general style) so only the cell
# contents change.

$dataRange = $ws.Range("BF2:BF31")
$dataRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "6-30-2007-08") {
        $cell.Value2 = "2008-06-30"
    }
}

$dataRange.ClearFormats()
